# Femacal de La Calera - Pepino ensalada: add the latest weekly price record.
#
# A new observation (week of 44785) is inserted above the existing row 324,
# pushing every subsequent record down by one row (old row 324 -> new row
# 325, ..., old row 399 -> new row 400). All the descriptive columns for the
# new record (market, region, product, quality, unit, origin, classification)
# stay the same as the record that used to sit at row 324 - only the date and
# the volume/price figures change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 324:399 down to 325:400, leaving a blank row 324 to fill in.
$ws.Rows.Item(324).Insert()

$ws.Range("A324").Value = 3
$ws.Range("B324").Value = "Femacal de La Calera"
$ws.Range("C324").Value = "Coquimbo"
$ws.Range("D324").Value = 44785
$ws.Range("E324").Value = 5
$ws.Range("F324").Value = 100112043
$ws.Range("G324").Value = "Pepino ensalada"
$ws.Range("H324").Value = "Sin especificar"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 105
$ws.Range("K324").Value = 18000
$ws.Range("L324").Value = 18500
$ws.Range("M324").Value = 18262
$ws.Range("N324").Value = "$/caja 70 unidades"
$ws.Range("O324").Value = "Región de Arica y Parinacota"
$ws.Range("P324").Value = 261
$ws.Range("Q324").Value = 70
$ws.Range("R324").Value = "Hortaliza"
